# Auto-generated edit script for extraction_Q2_Q3_v0.1.1.xlsx
# Implements RQ2/RQ3 data-extraction updates for ~5 newly reviewed papers
# (rows 29, 49, 50, 56, 75, 86) across the RQ2.1 / RQ2.2 / RQ2.3 / RQ3 sheets,
# plus the row-height and cursor/selection bookkeeping that Excel records
# when a reviewer scrolls through the sheet while filling it in.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("RQ2.1")
$ws2 = $wb.Worksheets.Item("RQ2.2")
$ws3 = $wb.Worksheets.Item("RQ2.3")
$ws4 = $wb.Worksheets.Item("RQ3")

# ------------------------------------------------------------------
# 1) Cell content - newly extracted data for each paper (row)
# ------------------------------------------------------------------

# --- RQ2.1 ---
$ws1.Range("F29").Value = "Business process orchestrator; Trust management system"
$ws1.Range("G29").Value = "Cloud"
$ws1.Range("F49").Value = "Device Authentication; Public Key Infrastructure; Data registry; Business process orchestrator"
$ws1.Range("G49").Value = "Cloud"
$ws1.Range("F50").Value = "Device, Data, and Service Authorisation;Data registry; Business process orchestrator"
$ws1.Range("G50").Value = "Cloud"
$ws1.Range("F56").Value = "Device, Data, and Service Authorisation"
$ws1.Range("G56").Value = "Cloud-Full & Fog-LW"
$ws1.Range("F75").Value = "Trust management system`n`tTrust rating record"
$ws1.Range("G75").Value = "Fog"
$ws1.Range("F86").Value = "Trust management system"
$ws1.Range("G86").Value = "Fog"

# --- RQ2.2 ---
$ws2.Range("F29").Value = "Resource exchange records"
$ws2.Range("G29").Value = "N/A"
$ws2.Range("H29").Value = "N/A"
$ws2.Range("I29").Value = "N/A"
$ws2.Range("J29").Value = "N/A"
$ws2.Range("F49").Value = "Data placement records; Resource exchange records"
$ws2.Range("G49").Value = "Sensor Readings"
$ws2.Range("H49").Value = "N/A"
$ws2.Range("I49").Value = "Cryptographic keys generation; Device Authentication"
$ws2.Range("J49").Value = "N/A"
$ws2.Range("F50").Value = "Sensor Reading Hashes; Device descriptions; Resource exchange records"
$ws2.Range("G50").Value = "N/A"
$ws2.Range("H50").Value = "Contract between resource providers and consumers"
$ws2.Range("I50").Value = "N/A"
$ws2.Range("J50").Value = "N/A"
$ws2.Range("F56").Value = "Authorisation requests and responses"
$ws2.Range("G56").Value = "N/A"
$ws2.Range("H56").Value = "Authorisation mechanism"
$ws2.Range("I56").Value = "N/A"
$ws2.Range("J56").Value = "N/A"
$ws2.Range("F75").Value = "Trust ratings"
$ws2.Range("G75").Value = "N/A"
$ws2.Range("H75").Value = "N/A"
$ws2.Range("I75").Value = "N/A"
$ws2.Range("J75").Value = "N/A"
$ws2.Range("F86").Value = "Trust ratings"
$ws2.Range("G86").Value = "N/A"
$ws2.Range("H86").Value = "N/A"
$ws2.Range("I86").Value = "Reputation calculation"
$ws2.Range("J86").Value = "N/A"

# --- RQ2.3 ---
$ws3.Range("F29").Value = 2
$ws3.Range("G29").Value = "blockchain"
$ws3.Range("H29").Value = "UTXO"
$ws3.Range("I29").Value = "N/A"
$ws3.Range("J29").Value = "Proof-of-work"
$ws3.Range("K29").Value = "public"
$ws3.Range("L29").Value = "N/A"
$ws3.Range("M29").Value = "In-house BC system"
$ws3.Range("F49").Value = 1
$ws3.Range("G49").Value = "blockchain"
$ws3.Range("H49").Value = "UTXO"
$ws3.Range("I49").Value = "N/A"
$ws3.Range("J49").Value = "Proof-of-work"
$ws3.Range("K49").Value = "public"
$ws3.Range("L49").Value = "N/A"
$ws3.Range("M49").Value = "bitcoin"
$ws3.Range("F50").Value = 1
$ws3.Range("G50").Value = "blockchain"
$ws3.Range("H50").Value = "account"
$ws3.Range("I50").Value = "installed"
$ws3.Range("J50").Value = "Proof-of-work"
$ws3.Range("K50").Value = "consortium"
$ws3.Range("L50").Value = "N/A"
$ws3.Range("M50").Value = "MultiChain"
$ws3.Range("F56").Value = 1
$ws3.Range("G56").Value = "blockchain"
$ws3.Range("H56").Value = "account"
$ws3.Range("I56").Value = "on-chain"
$ws3.Range("J56").Value = "Proof-of-work"
$ws3.Range("K56").Value = "public"
$ws3.Range("L56").Value = "N/A"
$ws3.Range("M56").Value = "Ethereum"
$ws3.Range("F75").Value = 1
$ws3.Range("G75").Value = "blockchain"
$ws3.Range("H75").Value = "account"
$ws3.Range("I75").Value = "on-chain"
$ws3.Range("J75").Value = "Proof-of-work"
$ws3.Range("K75").Value = "Public"
$ws3.Range("L75").Value = "N/A"
$ws3.Range("M75").Value = "Ethereum"
$ws3.Range("F86").Value = 1
$ws3.Range("G86").Value = "blockchain"
$ws3.Range("H86").Value = "UTXO"
$ws3.Range("I86").Value = "N/A"
$ws3.Range("J86").Value = "Joint PoW PoS"
$ws3.Range("K86").Value = "Public"
$ws3.Range("L86").Value = "N/A"
$ws3.Range("M86").Value = "in-house BC system"

# --- RQ3 ---
$ws4.Range("F29").Value = "N/A"
$ws4.Range("G29").Value = "N/A"
$ws4.Range("F49").Value = "N/A"
$ws4.Range("G49").Value = "N/A"
$ws4.Range("F56").Value = "N/A"
$ws4.Range("G56").Value = "N/A"
$ws4.Range("F75").Value = "N/A"
$ws4.Range("G75").Value = "N/A"
$ws4.Range("F86").Value = "N/A"
$ws4.Range("G86").Value = "N/A"

# ------------------------------------------------------------------
# 2) Row-height tweaks (rows whose wrapped text now needs more room)
# ------------------------------------------------------------------
$ws1.Rows.Item(49).RowHeight = 64
$ws2.Rows.Item(50).RowHeight = 80

# ------------------------------------------------------------------
# 3) Sheet scroll position / active-cell selection left by the reviewer.
#    Select() on the sheet last-activated (RQ3) must run last so the
#    workbook's active tab / tabSelected flag ends up unchanged.
# ------------------------------------------------------------------
$ws1.Range("G30").Select()
$ws2.Range("G30").Select()
$ws3.Range("H30").Select()
$ws4.Range("F30").Select()

